# Couch_Party / Prio.xlsx update ("mise a jour excel")
# Mark "Lobby" (row 16) and "Ecran Chargement" (row 18) as done:
#  - value goes from 0 (0%) to 1 (100%)
#  - conditional-style fill goes from red ("todo") to green ("done"),
#    matching the style already used by the other completed rows (e.g. B15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an already "done" cell (B15 - "Ecran D'accueil", 100%) as the style
# donor so the red->green fill/style swap exactly matches how Excel itself
# flips these status cells.
$doneStyleDonor = $ws.Range("B15")

$targets = @("B16", "B18")
foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    $doneStyleDonor.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats: bring over the green fill/number format
    $cell.Value = 1             # 100% done
}

$excel.CutCopyMode = 0

# Leave the selection where the author ended up after this edit.
$ws.Range("B18").Select() | Out-Null
